$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.797.08"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "2.120.69"
$ws.Range("E3").Value = "  +10.49%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.670"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.50%  "

$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("E11").Value = "  -3.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "2.428.71"
$ws.Range("E13").Value = "  +10.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.849"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.79%  "

$ws.Range("D16").Value = "2.119.27"
$ws.Range("E16").Value = "  +10.49%  "

$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "36.711.35"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("E21").Value = "  +2.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.98%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.48%  "

$ws.Range("E28").Value = "  +5.33%  "

$ws.Range("E29").Value = "  -7.48%  "

$ws.Range("E30").Value = "  -4.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +53.68%  "

$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0963"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.40%  "

$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("E35").Value = "  +18.03%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.35%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.34%  "

$ws.Range("E40").Value = "  -8.01%  "

$ws.Range("E41").Value = "  +7.09%  "

$ws.Range("E42").Value = "  -1.85%  "

$ws.Range("E43").Value = "  -6.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.86%  "

$ws.Range("D46").Value = "1.359.66"
$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.69%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0841"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.21%  "

$ws.Range("D49").Value = "2.319.12"
$ws.Range("E49").Value = "  +10.41%  "

$ws.Range("E50").Value = "  -3.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "

